$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
# Column A header: "name of file" -> "fileName"
$ws.Range("A1").Value = "fileName"
# Column C header: "price"
$ws.Range("C1").Value = "price"
# Column D header: "date"
$ws.Range("D1").Value = "date"

# Copy style of existing header (A1/B1) onto the new header cells C1, then
# give C1/D1 the same font/alignment/border treatment as A1/B1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update column B values: actg-2026 -> actg-2024 ---
$ws.Range("B2:B11").Value = "actg-2024"

# --- Add price column C ---
$prices = @(100.1234, 101.9999, 110.1234, 120.3579, 120.1234, 130.3579, 130.1234, 140.3579, 140.1234, 150.3579)
for ($i = 0; $i -lt $prices.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $prices[$i]
}

# --- Add date column D ---
$dates = @("01/01/2025","01/02/2026","01/03/2026","01/04/2026","01/05/2026","01/06/2026","01/07/2026","01/08/2026","01/09/2026","01/10/2026")
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dates[$i]
}
$ws.Range("D2:D11").NumberFormat = "mm-dd-yy"

# --- Borders on header row: thin left/right on C1:D1 like A1:B1 ---
# (Already copied via PasteSpecial above)

# --- Column widths (approximate autofit) ---
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

# --- Select A2 as the active cell (matches final selection state) ---
$ws.Range("A2").Select() | Out-Null

$ws.Range("A1:D11").Value
